$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old second header row (row 2). This shifts the existing data
# rows (old rows 3-11) up into rows 2-10, keeping their values/styles intact.
$ws.Rows(2).Delete()

# Rebuild row 1 as a single consolidated header row.
$ws.Range("A1:K1").ClearContents()
$ws.Range("A1:E1").ClearFormats()

$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Apply Arial 9 / General formatting to the new F1:K1 header cells. Adding a
# temporary named style lets us set the font without Excel also stamping an
# explicit "apply number format" flag onto the resulting cell style, and we
# remove the named style again afterwards so only the cellXfs entry remains.
$tmpStyle = $wb.Styles.Add("TmpHeaderStyle")
$tmpStyle.Font.Name = "Arial"
$tmpStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "TmpHeaderStyle"
$wb.Styles("TmpHeaderStyle").Delete()

$ws.Range("A2:K2").Select()
